# Fixing errors in example upload files.

$wb = $excel.ActiveWorkbook

# Remember which sheet/tab was active so we can restore it at the end -
# selecting ranges on other sheets below would otherwise change the
# workbook's active tab as a side effect.
$originalActiveSheetName = $wb.ActiveSheet.Name

# ---------------------------------------------------------------------------
# "Service Contacts" sheet (xl/worksheets/sheet5.xml):
#   - add a custom width for column A
#   - move the selection to D3
# ---------------------------------------------------------------------------
$wsServiceContacts = $wb.Worksheets.Item("Service Contacts")
$wsServiceContacts.Columns.Item(1).ColumnWidth = 13.666666666666666
$wsServiceContacts.Range("D3").Select()

# ---------------------------------------------------------------------------
# "Practitioners" sheet (xl/worksheets/sheet9.xml):
#   - add custom widths for columns A, C and F
#   - add a new data row (row 6)
#   - move the selection to G1:G1048576 (whole column G)
# ---------------------------------------------------------------------------
$wsPractitioners = $wb.Worksheets.Item("Practitioners")
$wsPractitioners.Columns.Item(1).ColumnWidth = 13.833333333333332
$wsPractitioners.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPractitioners.Columns.Item(6).ColumnWidth = 12.0

$wsPractitioners.Range("A6").Value = "PHN999:NFP02"
$wsPractitioners.Range("B6").Value = "P01"
$wsPractitioners.Range("C6").Value = 8
$wsPractitioners.Range("D6").Value = 1
$wsPractitioners.Range("E6").Value = 1973
$wsPractitioners.Range("F6").Value = 2
$wsPractitioners.Range("G6").Value = 1
$wsPractitioners.Range("H6").Value = 1
$wsPractitioners.Range("I6").Value = "tag1"

$wsPractitioners.Range("G1:G1048576").Select()

# ---------------------------------------------------------------------------
# Restore the original active sheet/tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($originalActiveSheetName).Activate()
